$d = $word.ActiveDocument

# Update the date line at the top of the document.
$d.Content.Find.Execute("2025-11-26 Wednesday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-11-27 Thursday", 2) | Out-Null

# Update each answer cell in the practice table. Several cells share the
# same original text (e.g. two cells both contain "42÷4=10, 2"), so a
# document-wide Find/Replace would clobber the wrong cell. Instead we
# address each cell directly via the Tables collection and assign its new
# value straight onto the cells Range, which leaves the run formatting
# (font/size) untouched and only rewrites that single cells text.
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "61÷9=6, 7"
$t.Cell(1, 2).Range.Text = "27÷3=9, 0"
$t.Cell(1, 3).Range.Text = "53÷7=7, 4"
$t.Cell(1, 4).Range.Text = "68÷8=8, 4"
$t.Cell(1, 5).Range.Text = "56÷4=14, 0"
$t.Cell(5, 1).Range.Text = "64÷5=12, 4"
$t.Cell(5, 2).Range.Text = "49÷6=8, 1"
$t.Cell(5, 3).Range.Text = "99÷7=14, 1"
$t.Cell(5, 4).Range.Text = "99÷3=33, 0"
$t.Cell(5, 5).Range.Text = "29÷8=3, 5"
$t.Cell(9, 1).Range.Text = "79÷3=26, 1"
$t.Cell(9, 2).Range.Text = "40÷8=5, 0"
$t.Cell(9, 3).Range.Text = "98÷6=16, 2"
$t.Cell(9, 4).Range.Text = "66÷6=11, 0"
$t.Cell(9, 5).Range.Text = "64÷8=8, 0"
$t.Cell(13, 1).Range.Text = "96÷5=19, 1"
$t.Cell(13, 2).Range.Text = "34÷2=17, 0"
$t.Cell(13, 3).Range.Text = "90÷3=30, 0"
$t.Cell(13, 4).Range.Text = "41÷6=6, 5"
$t.Cell(13, 5).Range.Text = "27÷5=5, 2"
$t.Cell(17, 1).Range.Text = "41÷6=6, 5"
$t.Cell(17, 2).Range.Text = "46÷6=7, 4"
$t.Cell(17, 3).Range.Text = "92÷4=23, 0"
$t.Cell(17, 4).Range.Text = "52÷4=13, 0"
$t.Cell(17, 5).Range.Text = "43÷4=10, 3"
